# Added post-irrad runs 12 and 13, loadings, photos
#
# The workbook keeps one worksheet per irradiation run; each new run's
# sheet is a copy of the previous run's sheet whose cells are formulas
# that mirror the prior sheet (e.g. ='Pre-irrad_11_141114'!A1). This
# script adds two new runs the same way:
#   Post-irrad_12_141117  (mirrors Pre-irrad_11_141114)
#   Post-irrad_13_141118  (mirrors Post-irrad_12_141117)

$wb = $excel.ActiveWorkbook

$dataRows = 17
$dataCols = 13  # A..M

# ---------------------------------------------------------------------
# Run 12: copy the last existing run sheet and drop it right after it.
# ---------------------------------------------------------------------
$prev1 = $wb.Worksheets.Item("Pre-irrad_11_141114")
$prev1.Copy($null, $prev1)
$run12 = $wb.Worksheets.Item($wb.Worksheets.Count)
$run12.Name = "Post-irrad_12_141117"

# The copy keeps formulas pointing at the sheet the *template* sheet
# referenced ('Pre-irrad_10_141113'); repoint them at the sheet this new
# tab actually mirrors ('Pre-irrad_11_141114').
for ($r = 1; $r -le $dataRows; $r++) {
  for ($c = 1; $c -le $dataCols; $c++) {
    $cell = $run12.Cells.Item($r, $c)
    $f = $cell.Formula
    if ($f -ne $null -and $f.StartsWith("=")) {
      $cell.Formula = $f.Replace("Pre-irrad_10_141113", "Pre-irrad_11_141114")
    }
  }
}

# ---------------------------------------------------------------------
# Run 13: same pattern, mirrors run 12.
# ---------------------------------------------------------------------
$run12.Copy($null, $run12)
$run13 = $wb.Worksheets.Item($wb.Worksheets.Count)
$run13.Name = "Post-irrad_13_141118"

for ($r = 1; $r -le $dataRows; $r++) {
  for ($c = 1; $c -le $dataCols; $c++) {
    $cell = $run13.Cells.Item($r, $c)
    $f = $cell.Formula
    if ($f -ne $null -and $f.StartsWith("=")) {
      $cell.Formula = $f.Replace("Pre-irrad_11_141114", "Post-irrad_12_141117")
    }
  }
}

# ---------------------------------------------------------------------
# Restore the editor's on-save selections / active sheet.
# ---------------------------------------------------------------------
$run12.Activate()
$run12.Range("B27").Select()

$run13.Activate()
$run13.Range("D34").Select()

$holderMatrix = $wb.Worksheets.Item("Holder Matrix, Irrad 141117")
$holderMatrix.Activate()
$holderMatrix.Range("A12").Select()

# Leave the newest run sheet as the active tab, matching the saved file.
$run13.Activate()
